$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing header cells in row 1 (strip numeric/typo suffixes, fix wording) ---
$ws.Range("A1").Value  = "Sender Cost Center"
$ws.Range("B1").Value  = "Personnel Number"
$ws.Range("C1").Value  = "Last name First name"
$ws.Range("D1").Value  = "Project Description"
$ws.Range("E1").Value  = "WBS Element"
$ws.Range("F1").Value  = "Network"
$ws.Range("G1").Value  = "Operation/Activity"
$ws.Range("H1").Value  = "Acct assgnt text"
$ws.Range("K1").Value  = "PROJ CODE"
$ws.Range("M1").Value  = "Raggr. X Struttura"
$ws.Range("O1").Value  = "Resp 2°liv"
$ws.Range("P1").Value  = "CDC bdg"
$ws.Range("Q1").Value  = "Resp 1°liv"
$ws.Range("R1").Value  = "Resp DT"
$ws.Range("S1").Value  = "TIPO REP"
$ws.Range("U1").Value  = "Bus Area_2"
$ws.Range("V1").Value  = "tipo progetto"
$ws.Range("W1").Value  = "Proj name"
$ws.Range("X1").Value  = "machine code"
$ws.Range("Y1").Value  = "product area"
$ws.Range("Z1").Value  = "Resp del Prodotto"
$ws.Range("AB1").Value = "Macchina"
$ws.Range("AF1").Value = "01.2025"
$ws.Range("AG1").Value = "02.2025"
$ws.Range("AH1").Value = "03.2025"
$ws.Range("AI1").Value = "04.2025"
$ws.Range("AJ1").Value = "05.2025"
$ws.Range("AK1").Value = "06.2025"
$ws.Range("AL1").Value = "07.2025"
$ws.Range("AM1").Value = "08.2025"
$ws.Range("AN1").Value = "09.2025"
$ws.Range("AO1").Value = "10.2025"
$ws.Range("AP1").Value = "11.2025"
$ws.Range("AQ1").Value = "12.2025"

# --- Add new column AZ with header "REP UT", matching style of column AY ---
$ws.Range("AZ1").Value = "REP UT"
$ws.Range("AZ1").Style = $ws.Range("AY1").Style

for ($r = 2; $r -le 13; $r++) {
    $src = $ws.Cells.Item($r, 51)   # AY
    $dst = $ws.Cells.Item($r, 52)   # AZ
    $dst.Style = $src.Style
}

# --- Column AZ width (matches the <col min="52" max="52" .../> added in the diff) ---
$ws.Columns.Item(52).ColumnWidth = 14.140625

# --- Sheet view: scroll to show the new columns, select the header row ---
$ws.Application.ActiveWindow.ScrollColumn = 30
$ws.Range("A1:AZ1").Select()
